# EV Charger Revenue Share by Recipient ISIC Code — split several aggregated
# ISIC categories on the "EVCRSbRIC" sheet into their finer-grained codes:
#   ISIC 05T06   -> ISIC 05, ISIC 06
#   ISIC 23      -> ISIC 231, ISIC 239
#   ISIC 24      -> ISIC 241, ISIC 242
#   ISIC 35T39   -> ISIC 351, ISIC 352T353, ISIC 36T39
# Each split inserts new column(s) immediately to the right of the original
# category column so the rest of the row shifts right, exactly like typing
# the new category name into a freshly inserted column in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVCRSbRIC")

# Insert the new (blank) columns working from right to left so that the
# column letters used below always refer to the original (pre-insert)
# positions for everything still to the left of the current insertion.
# Before state header layout (columns B..AL):
#   C  = ISIC 05T06
#   N  = ISIC 23
#   O  = ISIC 24
#   W  = ISIC 35T39

# ISIC 35T39 (col W) becomes 3 columns -> insert 2 new columns after W (at X)
$ws.Range("X1").EntireColumn.Insert()
$ws.Range("X1").EntireColumn.Insert()

# ISIC 24 (col O) becomes 2 columns -> insert 1 new column after O (at P)
$ws.Range("P1").EntireColumn.Insert()

# ISIC 23 (col N) becomes 2 columns -> insert 1 new column after N (at O)
$ws.Range("O1").EntireColumn.Insert()

# ISIC 05T06 (col C) becomes 2 columns -> insert 1 new column after C (at D)
$ws.Range("D1").EntireColumn.Insert()

# Set the header text for the split categories (column numbers in the final,
# post-insert layout).
$ws.Cells.Item(1,3).Value  = "ISIC 05"
$ws.Cells.Item(1,4).Value  = "ISIC 06"

$ws.Cells.Item(1,15).Value = "ISIC 231"
$ws.Cells.Item(1,16).Value = "ISIC 239"
$ws.Cells.Item(1,17).Value = "ISIC 241"
$ws.Cells.Item(1,18).Value = "ISIC 242"

$ws.Cells.Item(1,26).Value = "ISIC 351"
$ws.Cells.Item(1,27).Value = "ISIC 352T353"
$ws.Cells.Item(1,28).Value = "ISIC 36T39"

# The inserted columns leave row 2 ("Share of Costs by ISIC Code") blank in
# the new cells; every other non-formula cell on that row is a literal 0, so
# fill the newly created cells the same way.
$ws.Cells.Item(2,4).Value  = 0
$ws.Cells.Item(2,16).Value = 0
$ws.Cells.Item(2,18).Value = 0
$ws.Cells.Item(2,27).Value = 0
$ws.Cells.Item(2,28).Value = 0

Write-Output "EVCRSbRIC header/data columns updated"
